# Apply the edit described by the diff:
# 1. Add a "generic" word_type value in the new column J for the practice rows (2-5).
# 2. Add a new "stim details" table starting at row 27, with header row 28
#    (month / word_type / need_audio / need_image / word / count / find images)
#    followed by eight data rows (video x4, audio x4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column J values for the practice rows ---
$ws.Range("J2:J5").Value = "generic"

# --- New "stim details" block ---
$ws.Range("A27").Value = "stim details"

$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

$ws.Range("A29").Value = 6
$ws.Range("A30").Value = 6
$ws.Range("A31").Value = 7
$ws.Range("A32").Value = 7
$ws.Range("B29:B32").Value = "video"

$ws.Range("A33").Value = 6
$ws.Range("A34").Value = 6
$ws.Range("A35").Value = 7
$ws.Range("A36").Value = 7
$ws.Range("B33:B36").Value = "audio"

Write-Host "Edit applied."
